$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revise previously published monthly figures (small restatements) ---
$ws.Cells.Item(138, 10).Value = 3583212
$ws.Cells.Item(138, 13).Value = 365948
$ws.Cells.Item(138, 17).Value = 680512
$ws.Cells.Item(138, 23).Value = 3790727
$ws.Cells.Item(138, 24).Value = 473225

$ws.Cells.Item(139, 10).Value = 3446996
$ws.Cells.Item(139, 13).Value = 30604
$ws.Cells.Item(139, 17).Value = 110536
$ws.Cells.Item(139, 23).Value = 3882682
$ws.Cells.Item(139, 24).Value = -324838

$ws.Cells.Item(140, 10).Value = 5259010
$ws.Cells.Item(140, 13).Value = 519251
$ws.Cells.Item(140, 17).Value = -1012693
$ws.Cells.Item(140, 23).Value = 5903759
$ws.Cells.Item(140, 24).Value = -1656179

$ws.Cells.Item(141, 10).Value = 4819003
$ws.Cells.Item(141, 13).Value = 39214
$ws.Cells.Item(141, 17).Value = 1997897
$ws.Cells.Item(141, 23).Value = 5406220
$ws.Cells.Item(141, 24).Value = 1411009

$ws.Cells.Item(142, 10).Value = 5780784
$ws.Cells.Item(142, 13).Value = 31787
$ws.Cells.Item(142, 17).Value = -2790276
$ws.Cells.Item(142, 23).Value = 6338461
$ws.Cells.Item(142, 24).Value = -3347310

$ws.Cells.Item(143, 9).Value = 174623
$ws.Cells.Item(143, 10).Value = 5319569
$ws.Cells.Item(143, 13).Value = 22762
$ws.Cells.Item(143, 17).Value = -1586723
$ws.Cells.Item(143, 23).Value = 5991138
$ws.Cells.Item(143, 24).Value = -2258097

$ws.Cells.Item(144, 10).Value = 7004520
$ws.Cells.Item(144, 13).Value = 335949
$ws.Cells.Item(144, 17).Value = -2864822
$ws.Cells.Item(144, 23).Value = 7574645
$ws.Cells.Item(144, 24).Value = -3434372

$ws.Cells.Item(145, 10).Value = 6093933
$ws.Cells.Item(145, 13).Value = 31397
$ws.Cells.Item(145, 17).Value = -1463760
$ws.Cells.Item(145, 23).Value = 6684377
$ws.Cells.Item(145, 24).Value = -2053700

# --- Append new month: 01-09-2021 (row 146) ---
$dateCell = $ws.Cells.Item(146, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "01-09-2021"
$dateCell.Style = "Normal"

$ws.Cells.Item(146, 2).Value = 5345874
$ws.Cells.Item(146, 3).Value = 4216499
$ws.Cells.Item(146, 4).Value = 280321
$ws.Cells.Item(146, 5).Value = 239920
$ws.Cells.Item(146, 6).Value = 5437
$ws.Cells.Item(146, 7).Value = 102034
$ws.Cells.Item(146, 8).Value = 75252
$ws.Cells.Item(146, 9).Value = 426411
$ws.Cells.Item(146, 10).Value = 7111100
$ws.Cells.Item(146, 11).Value = 1119042
$ws.Cells.Item(146, 12).Value = 396101
$ws.Cells.Item(146, 13).Value = 514494
$ws.Cells.Item(146, 14).Value = 4349669
$ws.Cells.Item(146, 15).Value = 720658
$ws.Cells.Item(146, 16).Value = 11134
$ws.Cells.Item(146, 17).Value = -1765226
$ws.Cells.Item(146, 18).Value = 663253
$ws.Cells.Item(146, 19).Value = 248
$ws.Cells.Item(146, 20).Value = 346882
$ws.Cells.Item(146, 21).Value = 316619
$ws.Cells.Item(146, 22).Value = 5346121
$ws.Cells.Item(146, 23).Value = 7774601
$ws.Cells.Item(146, 24).Value = -2428480
$ws.Cells.Item(146, 25).Value = -3
